# Refresh scraped "want-to-go" counts (and a couple of min-price values) across all 4 sheets.
# Values below come from a re-scrape of the same bilibili-show source data; this script only
# pokes the handful of cells that changed between the two scrapes (see commit diff).
$wb = $excel.ActiveWorkbook

# --- Sheet 1: 展览 ---
$ws = $wb.Worksheets.Item(1)
$ws.Range("F3").Value = 152  # was 151
$ws.Range("F4").Value = 816  # was 815
$ws.Range("F5").Value = 866  # was 865
$ws.Range("F6").Value = 703  # was 699
$ws.Range("F7").Value = 1259  # was 1258
$ws.Range("G7").Value = 60  # was 54
$ws.Range("F9").Value = 867  # was 861
$ws.Range("G9").Value = 45  # was 40.5
$ws.Range("F10").Value = 718  # was 715
$ws.Range("F13").Value = 383  # was 382
$ws.Range("F15").Value = 1025  # was 1019
$ws.Range("F16").Value = 11350  # was 11318
$ws.Range("F17").Value = 656  # was 654
$ws.Range("F18").Value = 55  # was 54
$ws.Range("F22").Value = 287  # was 286
$ws.Range("F23").Value = 1801  # was 1798
$ws.Range("F26").Value = 496  # was 495
$ws.Range("F29").Value = 300  # was 298
$ws.Range("F31").Value = 268  # was 267
$ws.Range("F32").Value = 81  # was 80
$ws.Range("F33").Value = 104  # was 103
$ws.Range("F35").Value = 187  # was 186
$ws.Range("F37").Value = 1201  # was 1196
# --- Sheet 2: 演出 ---
$ws = $wb.Worksheets.Item(2)
$ws.Range("F7").Value = 151  # was 150
$ws.Range("F10").Value = 249  # was 247
$ws.Range("G11").Value = "已售罄"  # was 0, now sold out
$ws.Range("F16").Value = 326  # was 325
$ws.Range("F21").Value = 3  # was 2
# --- Sheet 3: 本地生活 ---
$ws = $wb.Worksheets.Item(3)
$ws.Range("F2").Value = 838  # was 836
# --- Sheet 4: 全部类型 ---
$ws = $wb.Worksheets.Item(4)
$ws.Range("F2").Value = 838  # was 836
$ws.Range("F4").Value = 152  # was 151
$ws.Range("F5").Value = 816  # was 815
$ws.Range("F7").Value = 866  # was 865
$ws.Range("F8").Value = 703  # was 699
$ws.Range("F9").Value = 1259  # was 1258
$ws.Range("G9").Value = 60  # was 54
$ws.Range("F12").Value = 151  # was 150
$ws.Range("F13").Value = 867  # was 861
$ws.Range("G13").Value = 45  # was 40.5
$ws.Range("F14").Value = 718  # was 715
$ws.Range("F17").Value = 1025  # was 1019
$ws.Range("F18").Value = 11350  # was 11318
$ws.Range("F19").Value = 249  # was 247
$ws.Range("F20").Value = 656  # was 654
$ws.Range("F21").Value = 55  # was 54
$ws.Range("F23").Value = 288  # was 286
$ws.Range("F24").Value = 1801  # was 1798
$ws.Range("F26").Value = 496  # was 495
$ws.Range("F33").Value = 326  # was 325
$ws.Range("F34").Value = 300  # was 298
$ws.Range("F37").Value = 268  # was 267
$ws.Range("F38").Value = 81  # was 80
$ws.Range("F39").Value = 104  # was 103
$ws.Range("F42").Value = 187  # was 186
$ws.Range("F46").Value = 1201  # was 1196
